$wb = $excel.ActiveWorkbook

# Cell updates to apply on both the "展览" and "全部类型" sheets
$updates = @{
    "F3" = 79
    "F6" = 27
    "F7" = 64
    "F9" = 252
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
